# Apply weekly-refresh style update: rotate the data among rows 4-10,12,13
# (row 11 is left untouched) for the Haba / Mapocho Venta Directa sheet.
# For each row we set Fecha (D), Volumen (J), Precio mínimo (K),
# Precio máximo (L), Precio promedio ponderado (M), Origen (O) and
# Precio $/Kg (P) to the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
    5  = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    6  = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    7  = @{ D = 44453; J = 55; K = 14000; L = 15000; M = 14455; O = "Provincia de Limarí"; P = 578 }
    8  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
    9  = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    10 = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    12 = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    13 = @{ D = 44376; J = 15; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
}
